$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was inserted as row 14, pushing every
# existing record from the old row 14 down through the old row 62 down by
# one row (old row 62 ends up at row 63). Insert a blank row at 14 first so
# all subsequent rows shift down, then populate the new row with its data.
$ws.Rows.Item(14).Insert()

$ws.Range("A14").Value = 9
$ws.Range("B14").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C14").Value = 'Metropolitana'
$ws.Range("D14").Value = 44462
$ws.Range("E14").Value = 13
$ws.Range("F14").Value = 100112022
$ws.Range("G14").Value = 'Arveja Verde'
$ws.Range("H14").Value = 'Perfection'
$ws.Range("I14").Value = 'Primera'
$ws.Range("J14").Value = 18
$ws.Range("K14").Value = 29000
$ws.Range("L14").Value = 30000
$ws.Range("M14").Value = 29500
$ws.Range("N14").Value = '$/malla 25 kilos'
$ws.Range("O14").Value = 'Provincia de Huasco'
$ws.Range("P14").Value = 1180
$ws.Range("Q14").Value = 25
$ws.Range("R14").Value = 'Hortaliza'
